# "Final Updates from end-of-semester"
# Slide 3 ("How About Custom Pipeline Steps?"), Content Placeholder 2:
# the paragraph "Let's walk through the Himanshu Chandra's " + "Colab"
# (two runs, second marked err="1", plus a trailing endParaRPr) becomes a
# single clean run reading "Let's walk through Himanshu Chandra's Colab"
# (the word "the" is dropped) with no leftover endParaRPr.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)          # "Content Placeholder 2"
$tr = $sh.TextFrame.TextRange

# Locate the target paragraph (3rd paragraph in this text frame).
$targetPara = $tr.Paragraphs(3, 1)

# Sanity-check we have the right paragraph before touching it.
if ($targetPara.Text -like "Let*Himanshu Chandra*Colab*") {

    $apos = [char]0x2019   # curly apostrophe used throughout the deck
    $newText = "Let" + $apos + "s walk through Himanshu Chandra" + $apos + "s Colab"

    # Remove the whole paragraph (its run(s), pPr and endParaRPr go with it),
    # then rebuild it fresh immediately before the following paragraph so the
    # new paragraph picks up that paragraph's indent level (lvl="1") and a
    # plain, single run of text with no trailing endParaRPr.
    $targetPara.Delete()

    $nextPara = $tr.Paragraphs(3, 1)   # now "So much good stuff in here!"
    [void]$nextPara.InsertBefore($newText + "`r")
}
